$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_3_8_10", 0.1908480773235842, -5.370626633474223, 0.2383582481846868, -0.04781911824640628, 0.895492672920227, 1.32693088054657, 2.265351057052612, 1.768540382385254),
    @("model_3_8_12", 0.2354883899290166, -4.378469061929707, 0.2635466654692527, 0.03788932344492202, 0.8460890650749207, 1.120275497436523, 2.190433263778687, 1.623879075050354),
    @("model_3_8_11", 0.2439473481137765, -4.357824605184342, 0.2841279563396467, 0.05630546169512374, 0.8367274403572083, 1.115975618362427, 2.129218339920044, 1.592795610427856),
    @("model_3_8_13", 0.2669798007701862, -3.824226369893614, 0.2952071693408393, 0.1003548196452629, 0.8112373352050781, 1.004832863807678, 2.096265316009521, 1.518447875976562),
    @("model_3_8_14", 0.2724210062392776, -3.724204649415809, 0.2961485522386552, 0.1076697171542496, 0.8052155375480652, 0.9839995503425598, 2.093465566635132, 1.506101608276367),
    @("model_3_8_15", 0.2727459020257459, -3.678669592919775, 0.291611490644733, 0.1068824569137503, 0.8048559427261353, 0.9745149612426758, 2.10696005821228, 1.5074303150177),
    @("model_3_8_8", 0.2728525650282257, -3.662827849169005, 0.3331934738140522, 0.1423999338439906, 0.80473792552948, 0.9712153673171997, 1.983282685279846, 1.447482943534851),
    @("model_3_8_20", 0.2804438225772318, -3.31453555688983, 0.2798103824674558, 0.1208860620738398, 0.7963365912437439, 0.8986698985099792, 2.142060041427612, 1.483794569969177),
    @("model_3_8_19", 0.2821842076117518, -3.328947645765955, 0.2859202143187868, 0.1250114755122707, 0.7944105863571167, 0.901671826839447, 2.123887538909912, 1.476831555366516),
    @("model_3_8_21", 0.2865571798602513, -3.266596800128108, 0.2916547150039379, 0.1338402949767732, 0.7895709276199341, 0.8886849284172058, 2.106831550598145, 1.461930155754089),
    @("model_3_8_16", 0.2898309803412343, -3.262252133358104, 0.3016898018281183, 0.1424458356052615, 0.7859477996826172, 0.887779951095581, 2.07698392868042, 1.447405457496643),
    @("model_3_8_18", 0.2930134072535455, -3.050777075306376, 0.2926673344308377, 0.1487802378149672, 0.7824258804321289, 0.8437319993972778, 2.103819370269775, 1.436714053153992),
    @("model_3_8_22", 0.2930903140445161, -3.128360438395145, 0.2976389055138378, 0.1478332991192627, 0.7823407053947449, 0.8598917126655579, 2.089032649993896, 1.438312292098999),
    @("model_3_8_9", 0.2952509681982877, -2.609081402137864, 0.3061436077586851, 0.1888120375983362, 0.7799494862556458, 0.7517316341400146, 2.063736915588379, 1.369147300720215),
    @("model_3_8_17", 0.2955455922617645, -3.079243365979237, 0.3018603962823707, 0.1545440373503043, 0.7796234488487244, 0.8496611714363098, 2.076476573944092, 1.426985740661621),
    @("model_3_8_23", 0.2955993759694669, -3.114710845798736, 0.303051290960104, 0.1532142489740227, 0.7795639634132385, 0.8570486903190613, 2.072934627532959, 1.429230213165283),
    @("model_3_8_24", 0.3023975863066404, -2.882974775947871, 0.3023685598617565, 0.1677871359561088, 0.7720402479171753, 0.8087806105613708, 2.074965238571167, 1.404633641242981),
    @("model_3_8_3", 0.3429506295188531, -0.9317797868472495, 0.3628248904556325, 0.3454001804179606, 0.7271600365638733, 0.4023683369159698, 1.895149946212769, 1.10485303401947),
    @("model_3_8_4", 0.3501248391675685, -0.797415271208046, 0.3666090603167635, 0.3573161421307254, 0.7192202806472778, 0.3743816614151001, 1.883894681930542, 1.084740877151489),
    @("model_3_8_6", 0.3784913497999721, -0.6611401007012438, 0.4020321315112625, 0.3955948321387033, 0.6878269910812378, 0.3459970653057098, 1.778535842895508, 1.020133018493652),
    @("model_3_8_7", 0.3788672535019445, -0.658928622491235, 0.4000058244374012, 0.3940589936573231, 0.6874108910560608, 0.3455364406108856, 1.784562706947327, 1.022725224494934),
    @("model_3_8_0", 0.3952463933010676, -0.5501009022239085, 0.5440041537596411, 0.5205828878603393, 0.6692840456962585, 0.322868824005127, 1.356268405914307, 0.8091744780540466),
    @("model_3_8_2", 0.3968947703512651, 0.02949790205187164, 0.4563826056021152, 0.4857873378694068, 0.6674597859382629, 0.2021448314189911, 1.616881132125854, 0.8679034113883972),
    @("model_3_8_1", 0.4055680631786877, -0.075005731798792, 0.4968529603873086, 0.5125209519313862, 0.657861053943634, 0.2239117920398712, 1.49651026725769, 0.8227816820144653),
    @("model_3_8_5", 0.4155175326974087, 0.2927537394549828, 0.4466925345886276, 0.4949509303066659, 0.6468499302864075, 0.1473115533590317, 1.645702362060547, 0.8524368405342102),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
